$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original columns: A=Arquivo, B=Empresa, C=Terminal, D=Operador, E=Porta
# Remove column A (Arquivo) -> columns become A=Empresa, B=Terminal, C=Operador, D=Porta
$ws.Columns.Item(1).Delete()

# Insert two new blank columns at D (before Porta) for CNPJ and IP
# -> A=Empresa, B=Terminal, C=Operador, D=(new CNPJ), E=(new IP), F=Porta
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).Insert()

# Header row for new columns
$ws.Range("D1").Value = "CNPJ"
$ws.Range("E1").Value = "IP"

# Match header style (bold, bordered, centered) from the existing header style
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# CNPJ is numeric-looking text, force text format so leading structure / full digit
# string is preserved as a string rather than being parsed into a number.
$ws.Range("D2:D7").NumberFormat = "@"
$ws.Range("D2:D7").Value = "313333395000141"
$ws.Range("D8:D13").NumberFormat = "@"
$ws.Range("D8:D13").Value = "31371695000141"

# IP addresses
$ws.Range("E2:E7").Value = "555.55.555.53"
$ws.Range("E8:E13").Value = "172.27.221.53"

# Reset cell style of the new CNPJ/data cells back to the plain (unstyled) look
# used by the rest of the data rows, since NumberFormat assignment above tagged
# them with a distinct style index.
$ws.Range("C2").Copy()
$ws.Range("D2:D13").PasteSpecial(-4122)
$excel.CutCopyMode = 0
